$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Build value from 7821 to 7829
$ws.Range("C2").Value = 7829

# Update the active selection to C2
$ws.Range("C2").Select()
